$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (OpenAI Embeddings) - updated benchmark values (B5:G5); H5 AVERAGE formula recalculates automatically
$ws.Range("B5").Value = 74.294243934442804
$ws.Range("C5").Value = 67.508716647607301
$ws.Range("D5").Value = 85.731714618792694
$ws.Range("E5").Value = 68.029270579492305
$ws.Range("F5").Value = 71.938878273698506
$ws.Range("G5").Value = 60.822865220460301

# Row 12 (BGE-large) - updated benchmark values (B12:G12); H12 AVERAGE formula recalculates automatically
$ws.Range("B12").Value = 64.189865503518504
$ws.Range("C12").Value = 53.1955242383737
$ws.Range("D12").Value = 84.528207751115801
$ws.Range("E12").Value = 61.523248158503698
$ws.Range("F12").Value = 54.642694910491798
$ws.Range("G12").Value = 34.956792802720599

# Update the active cell selection on the sheet (from J10 to K19)
$ws.Range("K19").Select()
